$wb = $excel.ActiveWorkbook

# OFF sheet - Week 16 row (A3 = "R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 217
$wsOff.Range("C3").Value = 146
$wsOff.Range("D3").Value = 46
$wsOff.Range("E3").Value = 29

# DEF sheet - Week 16 row (A3 = "R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 244
$wsDef.Range("C3").Value = 181
$wsDef.Range("D3").Value = 69
$wsDef.Range("E3").Value = 32
